# Update crypto price/volume data (refreshed GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.189.45'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.217.17'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.84'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.29'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +7.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.18'
$ws.Range("E10").Value = '  +3.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.45'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0939'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.10'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.103'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.553.11'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.870'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.230.28'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.069.17'
$ws.Range("E19").Value = '  +2.15%  '
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.48'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.72'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").Value = '  -5.60%  '
$ws.Range("E28").Value = '  -3.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.70'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.39'
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.47'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.06'
$ws.Range("E33").Value = '  +10.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.122'
$ws.Range("E34").Value = '  +2.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0777'
$ws.Range("E35").Value = '  +6.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.36'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.10'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0315'
$ws.Range("E40").Value = '  +5.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.24'
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.69'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.02'
$ws.Range("E43").Value = '  -5.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.11'
$ws.Range("E44").Value = '  +3.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.56'
$ws.Range("E45").Value = '  -4.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.196'
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.58'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.101'
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("E51").Value = '  +1.73%  '
